$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "550350664@qq.com"
$ws.Range("B2").Value = "Cody"
$ws.Range("C2").Value = "哈哈哈"
$ws.Range("D2").Value = "你好啊"
$ws.Range("E2").Value = "宠物.jpg"
$ws.Range("G2").Value = 1

# Delete rows 3 and 4 (shrinks used range to A1:G2)
$ws.Range("A3:G4").Delete()
